$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 9424
$ws1.Range("F7").Value = 848
$ws1.Range("F10").Value = 1152
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202405/BQV7zeWg1716290459878.jpeg"
$ws1.Range("F15").Value = 427
$ws1.Range("F16").Value = 90
$ws1.Range("F18").Value = 1287

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 9424
$ws4.Range("F8").Value = 848
$ws4.Range("F11").Value = 1152
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202405/BQV7zeWg1716290459878.jpeg"
$ws4.Range("F16").Value = 427
$ws4.Range("F17").Value = 91
$ws4.Range("F19").Value = 1287
